# Update Neg_Change and Pos_Change sheets with latest filtered market data.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws1.Cells.Item(2,1).Value = "ONGC"
$ws1.Cells.Item(2,2).Value = 238.5
$ws1.Cells.Item(2,3).Value = 239.4
$ws1.Cells.Item(2,4).Value = 237.02
$ws1.Cells.Item(2,5).Value = 238.1
$ws1.Cells.Item(2,6).Value = 4252781
$ws1.Cells.Item(2,7).Value = 8868641
$ws1.Cells.Item(2,8).Value = -0.5204698217009799
$ws1.Cells.Item(2,9).Value = "ONGC"
$ws1.Cells.Item(3,1).Value = "ADANIGREEN"
$ws1.Cells.Item(3,2).Value = 1024
$ws1.Cells.Item(3,3).Value = 1042
$ws1.Cells.Item(3,4).Value = 1024
$ws1.Cells.Item(3,5).Value = 1038.2
$ws1.Cells.Item(3,6).Value = 2615096
$ws1.Cells.Item(3,7).Value = 5927384
$ws1.Cells.Item(3,8).Value = -0.558811104527731
$ws1.Cells.Item(3,9).Value = "ADANIGREEN"
$ws1.Cells.Item(4,1).Value = "MOTHERSON"
$ws1.Cells.Item(4,2).Value = 120.49
$ws1.Cells.Item(4,3).Value = 121.48
$ws1.Cells.Item(4,4).Value = 119.97
$ws1.Cells.Item(4,5).Value = 121.3
$ws1.Cells.Item(4,6).Value = 14876284
$ws1.Cells.Item(4,7).Value = 30012240
$ws1.Cells.Item(4,8).Value = -0.5043261016172068
$ws1.Cells.Item(4,9).Value = "MOTHERSON"
$ws1.Cells.Item(5,1).Value = "DIVISLAB"
$ws1.Cells.Item(5,2).Value = 6430
$ws1.Cells.Item(5,3).Value = 6519
$ws1.Cells.Item(5,4).Value = 6410
$ws1.Cells.Item(5,5).Value = 6419
$ws1.Cells.Item(5,6).Value = 128635
$ws1.Cells.Item(5,7).Value = 264850
$ws1.Cells.Item(5,8).Value = -0.5143099867849726
$ws1.Cells.Item(5,9).Value = "DIVISLAB"
$ws1.Cells.Item(6,1).Value = "ASTRAL"
$ws1.Cells.Item(6,2).Value = 1410
$ws1.Cells.Item(6,3).Value = 1420.6
$ws1.Cells.Item(6,4).Value = 1404.7
$ws1.Cells.Item(6,5).Value = 1415.8
$ws1.Cells.Item(6,6).Value = 229029
$ws1.Cells.Item(6,7).Value = 510796
$ws1.Cells.Item(6,8).Value = -0.5516233486558234
$ws1.Cells.Item(6,9).Value = "ASTRAL"
$ws1.Cells.Item(7,1).Value = "TORNTPOWER"
$ws1.Cells.Item(7,2).Value = 1284
$ws1.Cells.Item(7,3).Value = 1296.5
$ws1.Cells.Item(7,4).Value = 1280.6
$ws1.Cells.Item(7,5).Value = 1289.8
$ws1.Cells.Item(7,6).Value = 143211
$ws1.Cells.Item(7,7).Value = 318739
$ws1.Cells.Item(7,8).Value = -0.5506950828107009
$ws1.Cells.Item(7,9).Value = "TORNTPOWER"
$ws1.Cells.Item(8,1).Value = "BANKINDIA"
$ws1.Cells.Item(8,2).Value = 141
$ws1.Cells.Item(8,3).Value = 141.57
$ws1.Cells.Item(8,4).Value = 139.83
$ws1.Cells.Item(8,5).Value = 141.2
$ws1.Cells.Item(8,6).Value = 3210687
$ws1.Cells.Item(8,7).Value = 7263380
$ws1.Cells.Item(8,8).Value = -0.5579624086857634
$ws1.Cells.Item(8,9).Value = "BANKINDIA"
$ws1.Cells.Item(9,1).Value = "CESC"
$ws1.Cells.Item(9,2).Value = 169.5
$ws1.Cells.Item(9,3).Value = 170.48
$ws1.Cells.Item(9,4).Value = 168.96
$ws1.Cells.Item(9,5).Value = 170.25
$ws1.Cells.Item(9,6).Value = 428230
$ws1.Cells.Item(9,7).Value = 877120
$ws1.Cells.Item(9,8).Value = -0.5117771798613645
$ws1.Cells.Item(9,9).Value = "CESC"
$ws1.Cells.Item(10,1).Value = "CDSL"
$ws1.Cells.Item(10,2).Value = 1531
$ws1.Cells.Item(10,3).Value = 1537.4
$ws1.Cells.Item(10,4).Value = 1518
$ws1.Cells.Item(10,5).Value = 1524
$ws1.Cells.Item(10,6).Value = 907070
$ws1.Cells.Item(10,7).Value = 1810907
$ws1.Cells.Item(10,8).Value = -0.4991073533870044
$ws1.Cells.Item(10,9).Value = "CDSL"
$ws1.Cells.Item(11,1).Value = "KFINTECH"
$ws1.Cells.Item(11,2).Value = 1080.9
$ws1.Cells.Item(11,3).Value = 1081.5
$ws1.Cells.Item(11,4).Value = 1045.6
$ws1.Cells.Item(11,5).Value = 1050.8
$ws1.Cells.Item(11,6).Value = 837476
$ws1.Cells.Item(11,7).Value = 1920257
$ws1.Cells.Item(11,8).Value = -0.5638729607547323
$ws1.Cells.Item(11,9).Value = "KFINTECH"

$ws2 = $wb.Worksheets.Item("Pos_Change")
$ws2.Cells.Item(2,1).Value = "MARUTI"
$ws2.Cells.Item(2,2).Value = 16248
$ws2.Cells.Item(2,3).Value = 16536
$ws2.Cells.Item(2,4).Value = 16247
$ws2.Cells.Item(2,5).Value = 16515
$ws2.Cells.Item(2,6).Value = 404711
$ws2.Cells.Item(2,7).Value = 268701
$ws2.Cells.Item(2,8).Value = 0.5061760097655015
$ws2.Cells.Item(2,9).Value = "MARUTI"
$ws2.Cells.Item(3,1).Value = "SHRIRAMFIN"
$ws2.Cells.Item(3,2).Value = 852
$ws2.Cells.Item(3,3).Value = 854.75
$ws2.Cells.Item(3,4).Value = 841.5
$ws2.Cells.Item(3,5).Value = 849
$ws2.Cells.Item(3,6).Value = 6036788
$ws2.Cells.Item(3,7).Value = 4213992
$ws2.Cells.Item(3,8).Value = 0.4325580115007338
$ws2.Cells.Item(3,9).Value = "SHRIRAMFIN"
$ws2.Cells.Item(4,1).Value = "JINDALSTEL"
$ws2.Cells.Item(4,2).Value = 1019
$ws2.Cells.Item(4,3).Value = 1032.5
$ws2.Cells.Item(4,4).Value = 1017.9
$ws2.Cells.Item(4,5).Value = 1029
$ws2.Cells.Item(4,6).Value = 846932
$ws2.Cells.Item(4,7).Value = 576973
$ws2.Cells.Item(4,8).Value = 0.4678884453865259
$ws2.Cells.Item(4,9).Value = "JINDALSTEL"
$ws2.Cells.Item(5,1).Value = "IOC"
$ws2.Cells.Item(5,2).Value = 162.35
$ws2.Cells.Item(5,3).Value = 164.9
$ws2.Cells.Item(5,4).Value = 162.06
$ws2.Cells.Item(5,5).Value = 163.85
$ws2.Cells.Item(5,6).Value = 13492656
$ws2.Cells.Item(5,7).Value = 8829634
$ws2.Cells.Item(5,8).Value = 0.5281104516902966
$ws2.Cells.Item(5,9).Value = "IOC"
$ws2.Cells.Item(6,1).Value = "CHOLAFIN"
$ws2.Cells.Item(6,2).Value = 1719.4
$ws2.Cells.Item(6,3).Value = 1740.4
$ws2.Cells.Item(6,4).Value = 1703.8
$ws2.Cells.Item(6,5).Value = 1734
$ws2.Cells.Item(6,6).Value = 2309890
$ws2.Cells.Item(6,7).Value = 1613630
$ws2.Cells.Item(6,8).Value = 0.4314867720605096
$ws2.Cells.Item(6,9).Value = "CHOLAFIN"
$ws2.Cells.Item(7,1).Value = "DLF"
$ws2.Cells.Item(7,2).Value = 697
$ws2.Cells.Item(7,3).Value = 705
$ws2.Cells.Item(7,4).Value = 691.7
$ws2.Cells.Item(7,5).Value = 698.3
$ws2.Cells.Item(7,6).Value = 2447290
$ws2.Cells.Item(7,7).Value = 1705673
$ws2.Cells.Item(7,8).Value = 0.4347943597629792
$ws2.Cells.Item(7,9).Value = "DLF"
$ws2.Cells.Item(8,1).Value = "BOSCHLTD"
$ws2.Cells.Item(8,2).Value = 36880
$ws2.Cells.Item(8,3).Value = 36950
$ws2.Cells.Item(8,4).Value = 36285
$ws2.Cells.Item(8,5).Value = 36625
$ws2.Cells.Item(8,6).Value = 11065
$ws2.Cells.Item(8,7).Value = 7160
$ws2.Cells.Item(8,8).Value = 0.5453910614525139
$ws2.Cells.Item(8,9).Value = "BOSCHLTD"
$ws2.Cells.Item(9,1).Value = "MFSL"
$ws2.Cells.Item(9,2).Value = 1711
$ws2.Cells.Item(9,3).Value = 1724.9
$ws2.Cells.Item(9,4).Value = 1707.5
$ws2.Cells.Item(9,5).Value = 1716.5
$ws2.Cells.Item(9,6).Value = 547322
$ws2.Cells.Item(9,7).Value = 370735
$ws2.Cells.Item(9,8).Value = 0.4763159669305569
$ws2.Cells.Item(9,9).Value = "MFSL"
$ws2.Cells.Item(10,1).Value = "PATANJALI"
$ws2.Cells.Item(10,2).Value = 534
$ws2.Cells.Item(10,3).Value = 539.4
$ws2.Cells.Item(10,4).Value = 530
$ws2.Cells.Item(10,5).Value = 536.3
$ws2.Cells.Item(10,6).Value = 4415565
$ws2.Cells.Item(10,7).Value = 2804718
$ws2.Cells.Item(10,8).Value = 0.5743347459530691
$ws2.Cells.Item(10,9).Value = "PATANJALI"
$ws2.Cells.Item(11,1).Value = "HDFCAMC"
$ws2.Cells.Item(11,2).Value = 2676
$ws2.Cells.Item(11,3).Value = 2681.3
$ws2.Cells.Item(11,4).Value = 2653.6
$ws2.Cells.Item(11,5).Value = 2669.1
$ws2.Cells.Item(11,6).Value = 1304915
$ws2.Cells.Item(11,7).Value = 895280
$ws2.Cells.Item(11,8).Value = 0.4575495934232865
$ws2.Cells.Item(11,9).Value = "HDFCAMC"
$ws2.Cells.Item(12,1).Value = "NCC"
$ws2.Cells.Item(12,2).Value = 161
$ws2.Cells.Item(12,3).Value = 161.74
$ws2.Cells.Item(12,4).Value = 158.6
$ws2.Cells.Item(12,5).Value = 160.46
$ws2.Cells.Item(12,6).Value = 3509185
$ws2.Cells.Item(12,7).Value = 2315642
$ws2.Cells.Item(12,8).Value = 0.515426391471566
$ws2.Cells.Item(12,9).Value = "NCC"
$ws2.Cells.Item(13,1).Value = "LAURUSLABS"
$ws2.Cells.Item(13,2).Value = 1016.1
$ws2.Cells.Item(13,3).Value = 1027.5
$ws2.Cells.Item(13,4).Value = 1006.1
$ws2.Cells.Item(13,5).Value = 1010.1
$ws2.Cells.Item(13,6).Value = 1244218
$ws2.Cells.Item(13,7).Value = 818322
$ws2.Cells.Item(13,8).Value = 0.5204503850562493
$ws2.Cells.Item(13,9).Value = "LAURUSLABS"
